# In-place rubric modification
# Rewrites the "Historico" sheet so that:
#   - A new "Calificación" (rating) column is inserted after "Proyecto"
#   - A new trailing column (duplicate of "Proyecto") is appended
#   - The "Cyber" project row becomes the first data row
#   - "Prototipo finalizado" is renamed to "Producto"
#   - Two numeric ratings (3 and 3.2) are recorded for two of the projects

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Area",  "Categoria", "Proyecto",                                    "Calificación", "Lider",                     "Video",                  "Poster",               "Equipo"),
    @("Cyber", "Concepto",  "dsfomo´ghdsrg",                                3,              "Marlon Martínez",           "dgfgj sdlgondskgdsfg",   "df´pdsfkg´pdgksdg",    "dsfomo´ghdsrg"),
    @("Nano",  "Prototipo", "Titulo para mi proyecto de prubas Gerry",      "",             "Gerry Deustúa Hernández",   "BOLDBGOSDBGOSANGSIGNSOPGSDG", "BOLDBGOSDBGOSANGSIGNSOPGSDG", "Titulo para mi proyecto de prubas Gerry"),
    @("Nano",  "Concepto",  "Test89",                                       "",             "Mikel Edel",                "rrrrrrrrrr",              "aaaaaaa",              "Test89"),
    @("Nexus", "Concepto",  "Robot automata para automatizar automatas",    "",             "Gerry Deustúa Hernández",   "sdgasdgasdg",             "sadgsadg",             "Robot automata para automatizar automatas"),
    @("Nano",  "Producto",  "Proyecto de prueba",                           3.2,            "Marlon Martínez",           "link",                    "link",                 "Proyecto de prueba"),
    @("Nexus", "Producto",  "fsfa",                                         "",             "Marlon Martínez",           "fasfsa",                  "fsa",                  "fsfa")
)

# Copy the formatting of column G (Poster) into column H (the new trailing
# column) before writing values, so the new column inherits the header /
# body styles already used by the rest of the table.
$ws.Range("G1:G7").Copy() | Out-Null
$ws.Range("H1:H7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth
